$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (Name, Position, Team), per the target diff.
$data = @(
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("Jalen Williams", "SG,SF,PF", "Oklahoma City Thunder"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
